# Update gh-pages to output generated at 456a3b4
# Applies numeric updates to the "想去人数" (F) column across sheets
# 展览 (Exhibition), 演出 (Performance), and 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 814
$wsExhibition.Range("F10").Value = 746
$wsExhibition.Range("F15").Value = 1226
$wsExhibition.Range("F25").Value = 1101

$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F9").Value = 32

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F6").Value = 814
$wsAllTypes.Range("F13").Value = 746
$wsAllTypes.Range("F21").Value = 1226
$wsAllTypes.Range("F31").Value = 32
$wsAllTypes.Range("F38").Value = 1101
